$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $origStyle = $rng.Style
    $rng.Value = "'" + $value
    $rng.Style = $origStyle
}

# Row 17/18: coins swap places (WrappedEther <-> ShibaInu)
$ws.Range("B17").Value = "ShibaInu"
$ws.Range("C17").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
Set-TextValue "D17" "0.0000165"
$ws.Range("E17").Value = "  -1.55%  "

$ws.Range("B18").Value = "WrappedEther"
$ws.Range("C18").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextValue "D18" "3.310.17"
$ws.Range("E18").Value = "  -0.45%  "

# Price / Volume(1h) refresh for remaining rows
Set-TextValue "D2" "66.022.69"
$ws.Range("E2").Value = "  -0.84%  "
Set-TextValue "D3" "3.310.89"
$ws.Range("E3").Value = "  -0.77%  "
$ws.Range("E4").Value = "  -0.02%  "
Set-TextValue "D5" "585.46"
$ws.Range("E5").Value = "  +2.09%  "
Set-TextValue "D6" "182.43"
$ws.Range("E7").Value = "  +3.30%  "
$ws.Range("E8").Value = "  +0.02%  "
Set-TextValue "D9" "3.309.31"
$ws.Range("E9").Value = "  -0.83%  "
$ws.Range("E10").Value = "  -3.39%  "
Set-TextValue "D11" "6.82"
$ws.Range("E11").Value = "  +2.56%  "
$ws.Range("E12").Value = "  -0.66%  "
Set-TextValue "D13" "3.887.70"
$ws.Range("E13").Value = "  -0.71%  "
$ws.Range("E14").Value = "  -2.78%  "
Set-TextValue "D15" "66.100.99"
$ws.Range("E15").Value = "  -0.88%  "
Set-TextValue "D16" "26.11"
$ws.Range("E16").Value = "  -3.32%  "
Set-TextValue "D19" "424.28"
$ws.Range("E19").Value = "  -2.75%  "
$ws.Range("E20").Value = "  -2.62%  "
Set-TextValue "D21" "13.13"
$ws.Range("E21").Value = "  -3.09%  "
Set-TextValue "D22" "7.38"
$ws.Range("E22").Value = "  -2.73%  "
Set-TextValue "D23" "71.73"
$ws.Range("E23").Value = "  -2.53%  "
$ws.Range("E24").Value = "  +0.12%  "
$ws.Range("E25").Value = "  +0.22%  "
Set-TextValue "D26" "3.461.83"
$ws.Range("E26").Value = "  -0.78%  "
$ws.Range("E27").Value = "  -1.03%  "
Set-TextValue "D28" "0.201"
$ws.Range("E28").Value = "  +5.40%  "
$ws.Range("E29").Value = "  -3.84%  "
Set-TextValue "D30" "8.86"
$ws.Range("E30").Value = "  -1.78%  "
Set-TextValue "D31" "0.999"
$ws.Range("E31").Value = "  +0.04%  "
Set-TextValue "D32" "1.91"
$ws.Range("E32").Value = "  -2.58%  "
Set-TextValue "D33" "22.33"
$ws.Range("E33").Value = "  -2.15%  "
$ws.Range("E34").Value = "  +0.06%  "
$ws.Range("E35").Value = "  -2.55%  "
Set-TextValue "D36" "6.54"
$ws.Range("E36").Value = "  -3.13%  "
$ws.Range("E37").Value = "  -4.82%  "
Set-TextValue "D38" "160.35"
$ws.Range("E38").Value = "  -1.70%  "
$ws.Range("E39").Value = "  -3.62%  "
Set-TextValue "D40" "2.870.93"
$ws.Range("E41").Value = "  -0.77%  "
Set-TextValue "D42" "26.25"
$ws.Range("E42").Value = "  -4.22%  "
Set-TextValue "D43" "0.757"
$ws.Range("E43").Value = "  -4.86%  "
$ws.Range("E44").Value = "  -2.82%  "
Set-TextValue "D45" "39.85"
$ws.Range("E45").Value = "  -0.77%  "
Set-TextValue "D46" "0.0659"
$ws.Range("E46").Value = "  -1.11%  "
Set-TextValue "D47" "5.89"
$ws.Range("E47").Value = "  -4.88%  "
$ws.Range("E48").Value = "  -3.23%  "
Set-TextValue "D49" "313.04"
$ws.Range("E49").Value = "  -2.64%  "
Set-TextValue "D50" "23.03"
$ws.Range("E50").Value = "  -5.90%  "
Set-TextValue "D51" "0.0270"
$ws.Range("E51").Value = "  -0.99%  "
